$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row right after the existing row 424 (row 425 onward
# shifts down by one to make room) and populate it with the new
# Sandia / Perú quote.
$ws.Rows(425).Insert()

$ws.Cells.Item(425, 1).Value = 3
$ws.Cells.Item(425, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(425, 3).Value = "Coquimbo"
$ws.Cells.Item(425, 4).Value = 44889
$ws.Cells.Item(425, 5).Value = 5
$ws.Cells.Item(425, 6).Value = 100112028
$ws.Cells.Item(425, 7).Value = "Sandia"
$ws.Cells.Item(425, 8).Value = "Sin especificar"
$ws.Cells.Item(425, 9).Value = "Primera"
$ws.Cells.Item(425, 10).Value = 250
$ws.Cells.Item(425, 11).Value = 800
$ws.Cells.Item(425, 12).Value = 850
$ws.Cells.Item(425, 13).Value = 826
$ws.Cells.Item(425, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(425, 15).Value = "Perú"
$ws.Cells.Item(425, 16).Value = 826
$ws.Cells.Item(425, 17).Value = 1
$ws.Cells.Item(425, 18).Value = "Hortaliza"
